# Update "想去人数" (want-to-go count) figures that were refreshed in the
# gh-pages data regeneration (commit 456a3b4).
#
# Sheet "展览" (Exhibitions) and sheet "全部类型" (All types) both contain the
# same underlying rows, so the same F-column updates apply to each.

$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7206
$ws1.Range("F4").Value = 117
$ws1.Range("F5").Value = 169
$ws1.Range("F6").Value = 7
$ws1.Range("F7").Value = 88
$ws1.Range("F8").Value = 602

# ---- Sheet: 全部类型 ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 7206
$ws4.Range("F5").Value = 117
$ws4.Range("F6").Value = 169
$ws4.Range("F7").Value = 7
$ws4.Range("F9").Value = 88
$ws4.Range("F10").Value = 602
